# Auto-generated script to apply scheduled market-price refresh values
# to the Behemoth_Profits workbook (per sheet: ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H2").Value = 1089.5
$ws.Range("J2").Value = 1233.3334
$ws.Range("L2").Value = 1233.3334
$ws.Range("N2").Value = -1459.3334

$ws.Range("H33").Value = 1708.5834
$ws.Range("I33").Value = 1777.8572
$ws.Range("J33").Value = 1611.6
$ws.Range("K33").Value = 1777.8572
$ws.Range("L33").Value = 1611.6
$ws.Range("M33").Value = -1548.8572
$ws.Range("N33").Value = -2069.6

$ws.Range("H40").Value = 3321.111
$ws.Range("I40").Value = 3618.889
$ws.Range("K40").Value = 3618.889
$ws.Range("M40").Value = -3443.889

$ws.Range("H53").Value = 832.381
$ws.Range("I53").Value = 420.7
$ws.Range("J53").Value = 1206.6364
$ws.Range("K53").Value = 420.7
$ws.Range("L53").Value = 1206.6364
$ws.Range("M53").Value = 216.3
$ws.Range("N53").Value = -2480.6364

$ws.Range("H64").Value = 3899.8
$ws.Range("I64").Value = 3166.3333
$ws.Range("K64").Value = 3166.3333
$ws.Range("M64").Value = -2918.3333

$ws.Range("H67").Value = 3899.8
$ws.Range("I67").Value = 3166.3333
$ws.Range("K67").Value = 3166.3333
$ws.Range("M67").Value = -2308.3333

$ws.Range("H74").Value = 3734.1765
$ws.Range("I74").Value = 3537
$ws.Range("J74").Value = 4375
$ws.Range("K74").Value = 3537
$ws.Range("L74").Value = 4375
$ws.Range("M74").Value = -2601
$ws.Range("N74").Value = -6247

$ws.Range("H76").Value = 7658.3477
$ws.Range("I76").Value = 8595.5
$ws.Range("K76").Value = 8595.5
$ws.Range("M76").Value = -8280.5

$ws.Range("H77").Value = 3734.1765
$ws.Range("I77").Value = 3537
$ws.Range("J77").Value = 4375
$ws.Range("K77").Value = 17685
$ws.Range("L77").Value = 21875
$ws.Range("M77").Value = -13005
$ws.Range("N77").Value = -31235

$ws.Range("H79").Value = 7658.3477
$ws.Range("I79").Value = 8595.5
$ws.Range("K79").Value = 8595.5
$ws.Range("M79").Value = -7503.5

$ws.Range("H92").Value = 326.79166
$ws.Range("I92").Value = 265.4091
$ws.Range("K92").Value = 265.4091
$ws.Range("M92").Value = 982.5908999999999

$ws.Range("H99").Value = 565.2727
$ws.Range("I99").Value = 551.9
$ws.Range("K99").Value = 1655.7
$ws.Range("M99").Value = -157.6999999999998

$ws.Range("H137").Value = 184232.22
$ws.Range("I137").Value = 244684.3
$ws.Range("K137").Value = 734052.8999999999
$ws.Range("M137").Value = -731502.8999999999


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H3").Value = 498.6
$ws.Range("I3").Value = 498.6
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 498.6
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -383.6
$ws.Range("N3").ClearContents()

$ws.Range("H32").Value = 6855654
$ws.Range("I32").Value = 8774231
$ws.Range("K32").Value = 8774231
$ws.Range("M32").Value = -8773944

$ws.Range("H97").Value = 1164.3
$ws.Range("I97").Value = 1030.2963
$ws.Range("K97").Value = 1030.2963
$ws.Range("M97").Value = -534.2963

$ws.Range("H110").Value = 1263.3572
$ws.Range("I110").Value = 1184.591
$ws.Range("K110").Value = 1184.591
$ws.Range("M110").Value = 860.4090000000001

$ws.Range("H132").Value = 4931.8613
$ws.Range("I132").Value = 2470.3572
$ws.Range("K132").Value = 7411.071599999999
$ws.Range("M132").Value = -4881.071599999999


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H80").Value = 1997.8572
$ws.Range("I80").Value = 1796.8
$ws.Range("J80").Value = 2500.5
$ws.Range("K80").Value = 1796.8
$ws.Range("L80").Value = 2500.5
$ws.Range("M80").Value = -798.8
$ws.Range("N80").Value = -4496.5

$ws.Range("H83").Value = 1997.8572
$ws.Range("I83").Value = 1796.8
$ws.Range("J83").Value = 2500.5
$ws.Range("K83").Value = 8984
$ws.Range("L83").Value = 12502.5
$ws.Range("M83").Value = -3992
$ws.Range("N83").Value = -22486.5

$ws.Range("H134").Value = 236706.11
$ws.Range("I134").Value = 1459.8857
$ws.Range("K134").Value = 4379.6571
$ws.Range("M134").Value = -1844.6571


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H122").Value = 1821.0952
$ws.Range("I122").Value = 2007.5834
$ws.Range("J122").Value = 1572.4445
$ws.Range("K122").Value = 6022.7502
$ws.Range("L122").Value = 4717.333500000001
$ws.Range("M122").Value = -3572.7502
$ws.Range("N122").Value = -9617.333500000001


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H4").Value = 21728444
$ws.Range("J4").Value = 42674316
$ws.Range("L4").Value = 128022948
$ws.Range("N4").Value = -128023172

$ws.Range("H113").Value = 1478.1875
$ws.Range("I113").Value = 824.75
$ws.Range("J113").Value = 1696
$ws.Range("K113").Value = 2474.25
$ws.Range("L113").Value = 5088
$ws.Range("M113").Value = -304.25
$ws.Range("N113").Value = -9428

$ws.Range("H132").Value = 1990.1333
$ws.Range("J132").Value = 2360.5
$ws.Range("L132").Value = 21244.5
$ws.Range("N132").Value = -26304.5

$ws.Range("H134").Value = 8561.049999999999
$ws.Range("J134").Value = 13746.542
$ws.Range("L134").Value = 41239.626
$ws.Range("N134").Value = -51379.626


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H5").Value = 150000
$ws.Range("I5").Value = 150000
$ws.Range("K5").Value = 150000
$ws.Range("M5").Value = -149888

$ws.Range("H113").Value = 4433.2144
$ws.Range("I113").Value = 4006
$ws.Range("K113").Value = 4006
$ws.Range("M113").Value = -1836

$ws.Range("H126").Value = 3897.2354
$ws.Range("I126").Value = 3423
$ws.Range("J126").Value = 4766.6665
$ws.Range("K126").Value = 10269
$ws.Range("L126").Value = 14299.9995
$ws.Range("M126").Value = -7799
$ws.Range("N126").Value = -19239.9995


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H16").Value = 1116
$ws.Range("I16").Value = 1116
$ws.Range("K16").Value = 1116
$ws.Range("M16").Value = -946

$ws.Range("H22").Value = 3091.6191
$ws.Range("J22").Value = 3394.2
$ws.Range("L22").Value = 3394.2
$ws.Range("N22").Value = -3984.2

$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("M24").ClearContents()
$ws.Range("N24").ClearContents()

$ws.Range("H27").Value = 3091.6191
$ws.Range("J27").Value = 3394.2
$ws.Range("L27").Value = 3394.2
$ws.Range("N27").Value = -3608.2

$ws.Range("H46").Value = 2765.7222
$ws.Range("I46").Value = 2209.5557
$ws.Range("K46").Value = 2209.5557
$ws.Range("M46").Value = -2021.5557

$ws.Range("H82").Value = 2842.1428
$ws.Range("J82").Value = 3001.5
$ws.Range("L82").Value = 3001.5
$ws.Range("N82").Value = -3723.5

$ws.Range("H85").Value = 2842.1428
$ws.Range("J85").Value = 3001.5
$ws.Range("L85").Value = 3001.5
$ws.Range("N85").Value = -5497.5

$ws.Range("H100").Value = 3569.9
$ws.Range("J100").Value = 3712.375
$ws.Range("L100").Value = 3712.375
$ws.Range("N100").Value = -4794.375


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H100").Value = 1088.7
$ws.Range("I100").Value = 1186.8462
$ws.Range("K100").Value = 2373.6924
$ws.Range("M100").Value = -1832.6924

$ws.Range("H132").Value = 253473.78
$ws.Range("I132").Value = 3194.875
$ws.Range("K132").Value = 9584.625
$ws.Range("M132").Value = -7054.625

$ws.Range("H136").Value = 7652.2
$ws.Range("I136").Value = 7419.4546
$ws.Range("K136").Value = 22258.3638
$ws.Range("M136").Value = -19708.3638

